# Swap the quantity/value figures (Code, Rate-with-tax, Qty, Value) that were
# accidentally entered on the wrong one of two adjacent rows sharing the same
# item description. Columns A (Sl No) and C (Item description) stay put;
# columns B, D, E, F, G are exchanged between each paired row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "D", "E", "F", "G")

# Each tuple is a pair of row numbers whose B/D/E/F/G values must be swapped.
$rowPairs = @(
    @(136, 137),
    @(146, 148),
    @(163, 164),
    @(233, 234),
    @(277, 278),
    @(299, 300),
    @(311, 312),
    @(356, 357),
    @(465, 466),
    @(467, 468),
    @(476, 477),
    @(479, 480),
    @(487, 488),
    @(603, 604),
    @(608, 609),
    @(717, 718)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "${col}${r1}"
        $addr2 = "${col}${r2}"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
